$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 0.114932203144002
$ws.Range("D6").Value = 1.0036620248797701

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 3.1745900716739599

$ws.Range("C8").Value = 0.22298186126252098
$ws.Range("D8").Value = 2.0918264198030503

$ws.Range("D9").Value = 0.27962171670971997

$ws.Range("D10").Value = 0.28068715855350901

$ws.Range("D11").Value = 3.0320802067846002

$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0.981299552298457

$ws.Range("C13").Value = 0.25285042388255702
$ws.Range("D13").Value = 1.3658974121237

$ws.Range("C14").Value = 0.211139508051647
$ws.Range("D14").Value = 2.08593233020291

$ws.Range("D15").Value = -0.51599196186857299

$ws.Range("D16").Value = 0.24777772335540099

$ws.Range("D17").Value = 0.91926315099982703

$ws.Range("A2").Select()
